$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 20, shifting old row 20 (Toplam) down to row 21
$ws.Rows.Item(20).Insert()

# Copy the formatting from row 19 into the newly blank row 20
$ws.Range("A19:U19").Copy()
$ws.Range("A20:U20").PasteSpecial(-4122)  # xlPasteFormats

$ws.Range("A20").Value = "AYDIN"
$ws.Range("B20").Value = "SULTANHİSAR"
$ws.Range("C20").Value = 127776
$ws.Range("D20").Value = "7820458686"
$ws.Range("E20").Value = "SULTANHİSAR MAL MÜDÜRLÜĞÜ"
$ws.Range("F20").Value = "İlçe Milli Eğitim Müdürlüğü"
$ws.Range("G20").Value = "MEM"
